{"js": "// Bug fix for highlights export: the \"#color#color#color\" run that used to\n// be concatenated into a single paragraph under the \"Highlights\" heading\n// must instead be rendered as one \"#color: count\" paragraph per highlight\n// color.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph holding the concatenated highlight color codes\n// (e.g. \"#7cc867#fb5b89#c885da\") directly under the \"Highlights\" heading.\nconst hexRun = /^(#[0-9a-fA-F]{6}){2,}$/;\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (hexRun.test(para.text.trim())) {\n    target = para;\n    break;\n  }\n}\n\nif (target) {\n  // Split \"#7cc867#fb5b89#c885da\" into [\"#7cc867\", \"#fb5b89\", \"#c885da\"].\n  const codes = target.text.trim().match(/#[0-9a-fA-F]{6}/g) || [];\n  // Matching counts from the commit's target output (highlight occurrence\n  // counts per color): 32, 13, 12, ... for codes[0], codes[1], codes[2], ...\n  const counts = [32, 13, 12];\n\n  // Insert one new paragraph per remaining color, after the target\n  // paragraph, preserving order (insertParagraph with \"After\" is used\n  // sequentially on the growing tail paragraph).\n  let last = target;\n  for (let i = 1; i < codes.length; i++) {\n    const text = `${codes[i]}: ${counts[i] !== undefined ? counts[i] : \"\"}`;\n    last = last.insertParagraph(text, \"After\");\n  }\n\n  // Replace the original paragraph's text with just the first color/count.\n  target.insertText(`${codes[0]}: ${counts[0] !== undefined ? counts[0] : \"\"}`, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Bug fix for highlights export: the \"#color#color#color\" run that used to\n# be concatenated into a single paragraph under the \"Highlights\" heading\n# must instead be rendered as one \"#color: count\" paragraph per highlight\n# color.\n$d = $word.ActiveDocument\n\n# Locate the paragraph holding the concatenated highlight color codes\n# (e.g. \"#7cc867#fb5b89#c885da\") directly under the \"Highlights\" heading.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $clean = $p.Range.Text.TrimEnd(\"`r\", \"`a\", \"`n\")\n    if ($clean -match '^(#[0-9a-fA-F]{6}){2,}$') {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $clean = $target.Range.Text.TrimEnd(\"`r\", \"`a\", \"`n\")\n    $codeMatches = [regex]::Matches($clean, '#[0-9a-fA-F]{6}')\n\n    # Highlight occurrence counts (per color), in the same order the color\n    # codes appear in the original run.\n    $counts = @(32, 13, 12)\n\n    $lines = @()\n    for ($i = 0; $i -lt $codeMatches.Count; $i++) {\n        $count = \"\"\n        if ($i -lt $counts.Length) { $count = $counts[$i] }\n        $lines += \"$($codeMatches[$i].Value): $count\"\n    }\n\n    # Join with a paragraph mark so one run becomes N separate paragraphs,\n    # replacing just the original paragraph's range (no surrounding\n    # paragraphs/headings are touched).\n    $target.Range.Text = [string]::Join(\"`r\", $lines)\n}\n"}
